$d = $word.ActiveDocument

function Add-HyperlinkParagraphBefore($paraIndex, $url) {
    # Inserts a new paragraph (containing a single hyperlink run whose
    # visible text is $url) immediately before the paragraph at $paraIndex.
    $para = $d.Paragraphs.Item($paraIndex)
    $insertStart = $para.Range.Start
    $insertPoint = $d.Range($insertStart, $insertStart)
    $insertPoint.InsertBefore("$url`r")
    $textRange = $d.Range($insertStart, $insertStart + $url.Length)
    $d.Hyperlinks.Add($textRange, $url, $null, $null, $url) | Out-Null
}

function Add-HyperlinkParagraphAfter($paraIndex, $url) {
    # Inserts a new paragraph (containing a single hyperlink run whose
    # visible text is $url) immediately after the paragraph at $paraIndex.
    # (Implemented as inserting "before" the following paragraph, since
    # inserting text right at a paragraph's End boundary causes it to merge
    # into the following paragraph instead of forming its own.)
    Add-HyperlinkParagraphBefore ($paraIndex + 1) $url
}

# Locate the paragraph that currently holds both the "_GoBack" bookmark and
# the "Javascript Basics" run (paragraph Range.Text includes a trailing
# paragraph-mark character, so trim it before comparing).
$jsParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd("`r")
    if ($t -eq "Javascript Basics") {
        $jsParaIndex = $i
        break
    }
}

# 1. Insert a new paragraph with the d3noob.org hyperlink right before the
#    "Javascript Basics" paragraph.
Add-HyperlinkParagraphBefore $jsParaIndex "http://www.d3noob.org/2014/02/d3js-elements.html"

# 2. The "Javascript Basics" paragraph shifted down by one because of the
#    insertion above. Split it into two paragraphs: one holding just the
#    bookmark, and one holding just the "Javascript Basics" text, by
#    inserting a paragraph break right before the text.
$jsParaIndex = $jsParaIndex + 1
$jsPara = $d.Paragraphs.Item($jsParaIndex)
$jsTextStart = $jsPara.Range.Start
$splitPoint = $d.Range($jsTextStart, $jsTextStart)
$splitPoint.InsertAfter("`r")

# 3. The hyperlink paragraph that used to directly follow the "Javascript
#    Basics" paragraph is now two paragraphs further down (bookmark
#    paragraph + text paragraph were inserted in between). Append a new
#    paragraph with the w3schools hyperlink right after it.
$oldHyperlinkParaIndex = $jsParaIndex + 2
Add-HyperlinkParagraphAfter $oldHyperlinkParaIndex "http://www.w3schools.com/jsref/jsref_toprecision.asp"

# 4. Insert one more blank paragraph right after the newly added hyperlink
#    paragraph (before the two pre-existing trailing blank paragraphs), by
#    inserting a lone paragraph break right before the first pre-existing
#    trailing blank paragraph.
$newHyperlinkParaIndex = $oldHyperlinkParaIndex + 1
$trailingBlankParaIndex = $newHyperlinkParaIndex + 1
$trailingBlankPara = $d.Paragraphs.Item($trailingBlankParaIndex)
$trailingBlankPara.Range.InsertParagraphBefore()
